$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transaction Details")

$ws.Range("A2").Value = "TR20240730015409284"
$ws.Range("C2").Value = "30-07-2024"
$ws.Range("D2").Value = "01:54:09:284"

$ws.Range("A3").Value = "TR20240730015409392"
$ws.Range("C3").Value = "30-07-2024"
$ws.Range("D3").Value = "01:54:09:392"

$ws.Range("A4").Value = "TR20240730015409501"
$ws.Range("C4").Value = "30-07-2024"
$ws.Range("D4").Value = "01:54:09:501"

$ws.Range("A5").Value = "TR20240731015409611"
$ws.Range("C5").Value = "31-07-2024"
$ws.Range("D5").Value = "01:54:09:612"

$ws.Range("A6").Value = "TR20240731015409720"
$ws.Range("C6").Value = "31-07-2024"
$ws.Range("D6").Value = "01:54:09:720"
